$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Data Mining Proces", $true, $false, $false, $false, $false, `
              $true, 1, $false, "Data Mining & Modeling Proces", 2)
